$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4070.389
$ws.Range("J70").Value = 4304.7856
$ws.Range("L70").Value = 12914.3568
$ws.Range("N70").Value = -13454.3568
$ws.Range("H73").Value = 4070.389
$ws.Range("J73").Value = 4304.7856
$ws.Range("L73").Value = 12914.3568
$ws.Range("N73").Value = -14786.3568
$ws.Range("H98").Value = 8497.608
$ws.Range("I98").Value = 8545.237999999999
$ws.Range("J98").Value = 7997.5
$ws.Range("K98").Value = 8545.237999999999
$ws.Range("L98").Value = 7997.5
$ws.Range("M98").Value = -7047.237999999999
$ws.Range("N98").Value = -10993.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents() | Out-Null
$ws.Range("H115").Value = 509
$ws.Range("I115").Value = 486.25
$ws.Range("K115").Value = 1458.75
$ws.Range("M115").Value = 108.25
$ws.Range("H122").Value = 8497.608
$ws.Range("I122").Value = 8545.237999999999
$ws.Range("J122").Value = 7997.5
$ws.Range("K122").Value = 25635.714
$ws.Range("L122").Value = 23992.5
$ws.Range("M122").Value = -23185.714
$ws.Range("N122").Value = -28892.5
$ws.Range("H131").Value = 102908.4
$ws.Range("I131").Value = 128336.125
$ws.Range("K131").Value = 385008.375
$ws.Range("M131").Value = -379968.375
$ws.Range("H141").Value = 9300
$ws.Range("I141").Value = 8160
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 24480
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -19300
$ws.Range("N141").Value = -55360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4066.85
$ws.Range("I32").Value = 3966.288
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 3966.288
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -3679.288
$ws.Range("N32").Value = -10574
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents() | Out-Null
$ws.Range("H36").Value = 6507.25
$ws.Range("I36").Value = 5333.3335
$ws.Range("K36").Value = 5333.3335
$ws.Range("M36").Value = -4987.3335
$ws.Range("H97").Value = 483.7857
$ws.Range("I97").Value = 504.2
$ws.Range("J97").Value = 432.75
$ws.Range("K97").Value = 504.2
$ws.Range("L97").Value = 432.75
$ws.Range("M97").Value = -8.199999999999989
$ws.Range("N97").Value = -1424.75
$ws.Range("H122").Value = 2447.6667
$ws.Range("I122").Value = 2181.9
$ws.Range("K122").Value = 6545.700000000001
$ws.Range("M122").Value = -4095.700000000001
$ws.Range("H132").Value = 3061.653
$ws.Range("I132").Value = 2528.4883
$ws.Range("K132").Value = 7585.4649
$ws.Range("M132").Value = -5055.4649

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1193.25
$ws.Range("I105").Value = 1140.2
$ws.Range("K105").Value = 1140.2
$ws.Range("M105").Value = 606.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1475.1666
$ws.Range("I22").Value = 1475.1666
$ws.Range("K22").Value = 1475.1666
$ws.Range("M22").Value = -1125.1666
$ws.Range("H31").Value = 5087.273
$ws.Range("J31").Value = 6230.357
$ws.Range("L31").Value = 6230.357
$ws.Range("N31").Value = -6820.357
$ws.Range("H34").Value = 5087.273
$ws.Range("J34").Value = 6230.357
$ws.Range("L34").Value = 6230.357
$ws.Range("N34").Value = -6634.357
$ws.Range("H86").Value = 7621.6665
$ws.Range("I86").Value = 6474.25
$ws.Range("K86").Value = 6474.25
$ws.Range("M86").Value = -5351.25
$ws.Range("H89").Value = 7621.6665
$ws.Range("I89").Value = 6474.25
$ws.Range("K89").Value = 32371.25
$ws.Range("M89").Value = -26755.25
$ws.Range("H139").Value = 143145.5
$ws.Range("J139").Value = 143145.5
$ws.Range("L139").Value = 143145.5
$ws.Range("N139").Value = -153425.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 209.4
$ws.Range("I13").Value = 15.666667
$ws.Range("K13").Value = 47.000001
$ws.Range("M13").Value = 120.999999
$ws.Range("H34").Value = 498.75
$ws.Range("J34").Value = 565
$ws.Range("L34").Value = 1695
$ws.Range("N34").Value = -1863
$ws.Range("H64").Value = 3220
$ws.Range("J64").Value = 5066.6665
$ws.Range("L64").Value = 15199.9995
$ws.Range("N64").Value = -15739.9995
$ws.Range("H67").Value = 3220
$ws.Range("J67").Value = 5066.6665
$ws.Range("L67").Value = 15199.9995
$ws.Range("N67").Value = -17071.9995
$ws.Range("H137").Value = 14837.4
$ws.Range("J137").Value = 24219
$ws.Range("L137").Value = 72657
$ws.Range("N137").Value = -82857

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 983.4
$ws.Range("I13").Value = 206.85715
$ws.Range("J13").Value = 2795.3333
$ws.Range("K13").Value = 206.85715
$ws.Range("L13").Value = 2795.3333
$ws.Range("M13").Value = -67.85714999999999
$ws.Range("N13").Value = -3073.3333
$ws.Range("H26").Value = 49998
$ws.Range("I26").Value = 49998
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 49998
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -49718
$ws.Range("N26").ClearContents() | Out-Null
$ws.Range("H50").Value = 49998
$ws.Range("I50").Value = 49998
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 49998
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -49500
$ws.Range("N50").ClearContents() | Out-Null
$ws.Range("H126").Value = 2828.923
$ws.Range("I126").Value = 1825.1428
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 5475.428400000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -3005.428400000001
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 2376.7715
$ws.Range("I132").Value = 1996.2916
$ws.Range("J132").Value = 3206.9092
$ws.Range("K132").Value = 5988.8748
$ws.Range("L132").Value = 9620.7276
$ws.Range("M132").Value = -3458.8748
$ws.Range("N132").Value = -14680.7276

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5730.857
$ws.Range("I7").Value = 5635.8
$ws.Range("K7").Value = 5635.8
$ws.Range("M7").Value = -5523.8
$ws.Range("H22").Value = 1971.1428
$ws.Range("I22").Value = 1116
$ws.Range("K22").Value = 1116
$ws.Range("M22").Value = -821
$ws.Range("H27").Value = 1971.1428
$ws.Range("I27").Value = 1116
$ws.Range("K27").Value = 1116
$ws.Range("M27").Value = -1009
$ws.Range("H40").Value = 3200.7727
$ws.Range("I40").Value = 3021.3333
$ws.Range("K40").Value = 3021.3333
$ws.Range("M40").Value = -2885.3333
$ws.Range("H55").Value = 695.2
$ws.Range("I55").Value = 836.0909
$ws.Range("J55").Value = 523
$ws.Range("K55").Value = 836.0909
$ws.Range("L55").Value = 523
$ws.Range("M55").Value = -663.0909
$ws.Range("N55").Value = -869
$ws.Range("H61").Value = 170542.17
$ws.Range("I61").Value = 253125.75
$ws.Range("K61").Value = 253125.75
$ws.Range("M61").Value = -252923.75
$ws.Range("H74").Value = 38497.332
$ws.Range("I74").Value = 38497.332
$ws.Range("K74").Value = 38497.332
$ws.Range("M74").Value = -37499.332
$ws.Range("H77").Value = 38497.332
$ws.Range("I77").Value = 38497.332
$ws.Range("K77").Value = 115491.996
$ws.Range("M77").Value = -110499.996
$ws.Range("H113").Value = 170542.17
$ws.Range("I113").Value = 253125.75
$ws.Range("K113").Value = 253125.75
$ws.Range("M113").Value = -250955.75
$ws.Range("H122").Value = 2860.3333
$ws.Range("I122").Value = 1992.1666
$ws.Range("K122").Value = 5976.4998
$ws.Range("M122").Value = -3526.4998
$ws.Range("H126").Value = 5730.857
$ws.Range("I126").Value = 5635.8
$ws.Range("K126").Value = 16907.4
$ws.Range("M126").Value = -14437.4
$ws.Range("H132").Value = 9520.125
$ws.Range("I132").Value = 11105.823
$ws.Range("J132").Value = 5669.143
$ws.Range("K132").Value = 33317.469
$ws.Range("L132").Value = 17007.429
$ws.Range("M132").Value = -30787.469
$ws.Range("N132").Value = -22067.429

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4033.9167
$ws.Range("I81").Value = 1321.2858
$ws.Range("K81").Value = 2642.5716
$ws.Range("M81").Value = -1581.5716
$ws.Range("H84").Value = 4033.9167
$ws.Range("I84").Value = 1321.2858
$ws.Range("K84").Value = 13212.858
$ws.Range("M84").Value = -7908.858
$ws.Range("H107").Value = 1674.5
$ws.Range("I107").Value = 961.6667
$ws.Range("J107").Value = 2387.3333
$ws.Range("K107").Value = 2885.0001
$ws.Range("L107").Value = 7161.999899999999
$ws.Range("M107").Value = -965.0001000000002
$ws.Range("N107").Value = -11001.9999
$ws.Range("H113").Value = 417.5
$ws.Range("I113").Value = 418.8095
$ws.Range("J113").Value = 390
$ws.Range("K113").Value = 1256.4285
$ws.Range("L113").Value = 1170
$ws.Range("M113").Value = 913.5715
$ws.Range("N113").Value = -5510
$ws.Range("H136").Value = 6190.59
$ws.Range("I136").Value = 5307.8965
$ws.Range("J136").Value = 8750.4
$ws.Range("K136").Value = 15923.6895
$ws.Range("L136").Value = 26251.2
$ws.Range("M136").Value = -13373.6895
$ws.Range("N136").Value = -31351.2
$ws.Range("H139").Value = 85099.5
$ws.Range("I139").Value = 85099
$ws.Range("K139").Value = 85099
$ws.Range("M139").Value = -79959
